# MARCELINO_RAMOS.xlsx — structural update:
#  - delete the "Desarquivamentos Pendentes" sheet (no longer needed)
#  - rename "Paineis DARQ" -> "PAINEIS DARQ"
#  - rename "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO"

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

$wb.Worksheets.Item("Desarquivamentos Pendentes").Delete()

$wb.Worksheets.Item("Paineis DARQ").Name = "PAINEIS DARQ"
$wb.Worksheets.Item("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"
